$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update the Lacustrine-IUCNGET raw data path (B5) and drop its hyperlink
$ws.Range("B5").Hyperlinks.Delete()
$ws.Range("B5").Value = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\processing\NEAP_intermediate\Lakes_NEAP_20240808_NoOverlapWithALUM.shp"

# Update the ALUM_2010-IUCNGET raw data path (B7)
$ws.Range("B7").Value = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\inputs\raw\Land_use_of_Australia\ABARES_Land_use_of_Australia_2010_11_to_2020_21_prerelease3_20240809\NLUM_v7p3_ALUMV8_250m_2010_11_alb\NLUM_v7p3_ALUMV8_250m_2010_11_alb.tif"
